# Apply cryptocurrency price/volume updates to Sheet1
# Generated from the OOXML diff (92 cell changes across D and E columns, rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.215.57'
$ws.Range('E2').Value = '  -5.03%  '
$ws.Range('D3').Value = '3.000.60'
$ws.Range('E3').Value = '  -5.37%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.46'
$ws.Range('E5').Value = '  -4.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '125.62'
$ws.Range('E6').Value = '  -7.26%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '2.994.16'
$ws.Range('E8').Value = '  -5.62%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.503'
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.132'
$ws.Range('E10').Value = '  -7.28%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.06'
$ws.Range('E11').Value = '  -5.20%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.441'
$ws.Range('E12').Value = '  -2.87%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000221'
$ws.Range('E13').Value = '  -7.55%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.55'
$ws.Range('E14').Value = '  -5.73%  '
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = '3.500.26'
$ws.Range('E16').Value = '  -5.20%  '
$ws.Range('D17').Value = '3.001.24'
$ws.Range('E17').Value = '  -5.40%  '
$ws.Range('D18').Value = '60.238.57'
$ws.Range('E18').Value = '  -4.98%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.51'
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '431.11'
$ws.Range('E20').Value = '  -6.54%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.16'
$ws.Range('E21').Value = '  -5.62%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.673'
$ws.Range('E22').Value = '  -3.35%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.07'
$ws.Range('E23').Value = '  -7.68%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.89'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '79.56'
$ws.Range('E25').Value = '  -4.32%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.53'
$ws.Range('E28').Value = '  -5.72%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.96'
$ws.Range('E29').Value = '  -4.95%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.22'
$ws.Range('E30').Value = '  -6.63%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.15'
$ws.Range('E31').Value = '  -9.61%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '25.27'
$ws.Range('E32').Value = '  -7.16%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0950'
$ws.Range('E33').Value = '  -5.83%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.61'
$ws.Range('E34').Value = '  -4.64%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.933'
$ws.Range('E35').Value = '  -8.66%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '50.30'
$ws.Range('E36').Value = '  -2.15%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.03'
$ws.Range('E37').Value = '  -15.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.53'
$ws.Range('E38').Value = '  +5.04%  '
$ws.Range('D39').Value = '0.0₃0662'
$ws.Range('E39').Value = '  -10.54%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0358'
$ws.Range('E40').Value = '  -8.23%  '
$ws.Range('E41').Value = '  -4.19%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '372.19'
$ws.Range('E42').Value = '  -5.05%  '
$ws.Range('D43').Value = '2.676.95'
$ws.Range('E43').Value = '  -3.99%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.44'
$ws.Range('E44').Value = '  -7.10%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '121.61'
$ws.Range('E46').Value = '  -4.54%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.234'
$ws.Range('E47').Value = '  -6.66%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.99'
$ws.Range('E48').Value = '  -5.81%  '
$ws.Range('E49').Value = '  -3.57%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.35'
$ws.Range('E50').Value = '  -6.60%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.132'
$ws.Range('E51').Value = '  -2.24%  '
